# Apply trade #24 close update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# --- Summary sheet -----------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.3   # Current Capital
$wsSummary.Range("B4").Value = 0.3      # Total P&L $
$wsSummary.Range("B5").Value = 0.25     # Total P&L %
$wsSummary.Range("B6").Value = 24       # Total Trades
$wsSummary.Range("B7").Value = 9        # Winning Trades
$wsSummary.Range("B9").Value = 37.5     # Win Rate %

# --- Strategy Status sheet ----------------------------------------------
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 100.3   # Capital (MarketMaking)
$wsStrategy.Range("D4").Value = 24      # Trades
$wsStrategy.Range("E4").Value = 0.3     # P&L $
$wsStrategy.Range("F4").Value = 0.3     # P&L %
$wsStrategy.Range("G4").Value = 37.5    # Win Rate %

# --- New closed trade row, appended to both trade-log sheets -----------
$newTrade = @(24, "2026-02-17", "04:08:36", "MarketMaking", "UP", 0.35, 0.58, "CLOSED", 65.71429999999999, 0.23, 100.3, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 25

    # Date (B) and Time (C) columns look numeric to Excel's auto-detection,
    # so force them to stay text (matching the source's inline-string date
    # and time cells) before writing the values.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("C$row").NumberFormat = "@"

    for ($i = 0; $i -lt $newTrade.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($row, $col).Value = $newTrade[$i]
    }
}
